$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.997.95'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.913.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.87%  '

$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.53%  '

$ws.Range("E6").Value = '  -0.27%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4594'
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = '  -1.27%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07727'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9819'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.03%  '

$ws.Range("E11").Value = '  -2.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.907.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.06%  '

$ws.Range("E13").Value = '  -2.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.670'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07032'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.51%  '

$ws.Range("E16").Value = '  -0.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '83.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009470'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.96%  '

$ws.Range("E19").Value = '  -1.70%  '

$ws.Range("E20").Value = '  -0.40%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.976.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.61%  '

$ws.Range("E22").Value = '  -2.94%  '

$ws.Range("E23").Value = '  -1.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.70%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.671'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.95%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.856'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09305'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8680'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.082'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.254'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.85%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.024'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05732'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.155'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.58%  '

$ws.Range("E37").Value = '  -0.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02042'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.24%  '

$ws.Range("E39").Value = '  -2.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.407'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1755'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.857'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.337'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5193'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.30'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06872'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000002622'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.056'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.785'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.54'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.31%  '
